$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull/push data corrections to column F (dSF) -- recalculated mean values
$ws.Range("F16").Value = -1
$ws.Range("F19").Value = 3
$ws.Range("F25").Value = -5
$ws.Range("F31").Value = 3
$ws.Range("F32").Value = 4
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = -4
